$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 22, shifting existing rows 22-84 down to 23-85
$ws.Rows.Item(22).Insert()

# Fill in the new row 22 with the new data point
$ws.Cells.Item(22, 1).Value = 5
$ws.Cells.Item(22, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(22, 3).Value = "Maule"
$ws.Cells.Item(22, 4).Value = 44560
$ws.Cells.Item(22, 5).Value = 7
$ws.Cells.Item(22, 6).Value = 100112030
$ws.Cells.Item(22, 7).Value = "Poroto granado"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 300
$ws.Cells.Item(22, 11).Value = 23000
$ws.Cells.Item(22, 12).Value = 23000
$ws.Cells.Item(22, 13).Value = 23000
$ws.Cells.Item(22, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(22, 15).Value = "Región del Maule"
$ws.Cells.Item(22, 16).Value = 920
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"
